# New Submission Synced: 2026-02-08 19:21:38
# Target sheet: "JSS 3A" (this workbook's Google-Forms-synced results tab).
# 1) C6 ("Admission No" for Amina Abubakar Adam) was stored as text "47";
#    fix it to a genuine number 47.
# 2) Append the new form submission as row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# --- Fix C6: text "47" -> numeric 47 ---
$ws.Range("C6").Value = 47

# --- Append new submission row 7 ---
$ws.Range("A7").Value = "2026-02-08 19:21:37"
$ws.Range("B7").Value = "MUHAMMAD ABUBAKAR "

# C7 ("Admission No") is synced as text "5" (matches the other rows' source
# data, which are free-text form answers) rather than a number, so force
# text formatting before writing it.
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "5"

$ws.Range("D7").Value = 9
